$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.087680697441101
$ws.Range("B1").Value = 1.001729249954224
$ws.Range("C1").Value = 5.356661319732666
$ws.Range("D1").Value = 1.898563742637634
$ws.Range("E1").Value = 1.10181450843811
